$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (R) that mirrors the formatting of the existing
# "2020" column (Q) for rows 4-14, then fill in the 2021 values.
$srcFmt = $ws.Range("Q4:Q14")
$dstFmt = $ws.Range("R4:R14")
$srcFmt.Copy($dstFmt) | Out-Null

$ws.Cells.Item(4, 18).Value2 = 2021

$ws.Cells.Item(5, 18).Value2 = 5.8
$ws.Cells.Item(6, 18).Value2 = 4.7
$ws.Cells.Item(7, 18).Value2 = 1.6
$ws.Cells.Item(8, 18).Value2 = 12.9
$ws.Cells.Item(9, 18).Value2 = 10.199999999999999
$ws.Cells.Item(10, 18).Value2 = 4.2
$ws.Cells.Item(11, 18).Value2 = 3.3
$ws.Cells.Item(12, 18).Value2 = 15.2
$ws.Cells.Item(13, 18).Value2 = 2.4
$ws.Cells.Item(14, 18).Value2 = 0.6

# Move the active selection from T1 (whole column) to T9, matching the
# author's saved view state.
$ws.Range("T9").Select() | Out-Null
